$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '64.361.49'
$ws.Range('E2').Value = '  -0.03%  '

$ws.Range('D3').Value = '3.135.75'
$ws.Range('E3').Value = '  -0.74%  '

$ws.Range('E4').Value = '  -0.07%  '

$ws.Range('E5').Value = '  +0.21%  '

$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '143.29'
$ws.Range('E6').Value = '  -3.02%  '

$ws.Range('E7').Value = '  +0.07%  '

$ws.Range('D8').Value = '3.130.42'
$ws.Range('E8').Value = '  -0.86%  '

$ws.Range('E9').Value = '  +0.82%  '

$ws.Range('E10').Value = '  -0.86%  '

$ws.Range('E11').Value = '  -2.38%  '

$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.469'
$ws.Range('E12').Value = '  -1.68%  '

$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '0.0000254'
$ws.Range('E13').Value = '  +0.99%  '

$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '35.30'
$ws.Range('E14').Value = '  -1.00%  '

$ws.Range('D15').Value = '3.649.35'
$ws.Range('E15').Value = '  -0.73%  '

$ws.Range('D17').Value = '64.304.99'
$ws.Range('E17').Value = '  -0.07%  '

$ws.Range('D18').Value = '3.134.63'
$ws.Range('E18').Value = '  -0.72%  '

$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '6.84'
$ws.Range('E19').Value = '  -1.44%  '

$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '475.39'
$ws.Range('E20').Value = '  -1.35%  '

$ws.Range('E21').Value = '  +0.30%  '

$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '0.716'
$ws.Range('E22').Value = '  +0.79%  '

$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '7.78'
$ws.Range('E23').Value = '  -0.08%  '

$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '85.94'
$ws.Range('E24').Value = '  +2.62%  '

$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '13.51'
$ws.Range('E25').Value = '  -1.51%  '

$ws.Range('E26').Value = '  -0.01%  '

$ws.Range('E27').Value = '  -3.93%  '

$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '8.44'
$ws.Range('E28').Value = '  -0.42%  '

$ws.Range('E29').Value = '  +7.64%  '

$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '2.05'
$ws.Range('E30').Value = '  -6.55%  '

$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '0.114'
$ws.Range('E31').Value = '  +1.05%  '

$ws.Range('E32').Value = '  -0.03%  '

$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '26.60'
$ws.Range('E33').Value = '  +1.25%  '

$ws.Range('E34').Value = '  -3.62%  '

$ws.Range('E35').Value = '  +0.13%  '

$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '5.96'
$ws.Range('E36').Value = '  -0.53%  '

$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '52.65'

$ws.Range('E38').Value = '  +1.43%  '

$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '448.34'
$ws.Range('E39').Value = '  -2.21%  '

$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '2.98'
$ws.Range('E40').Value = '  +1.55%  '

$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.0393'
$ws.Range('E41').Value = '  -0.80%  '

$ws.Range('E42').Value = '  -0.46%  '

$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '8.30'
$ws.Range('E43').Value = '  -1.69%  '

$ws.Range('D44').Value = '2.877.49'
$ws.Range('E44').Value = '  +0.72%  '

$ws.Range('E45').Value = '  -1.67%  '

$ws.Range('E46').Value = '  -1.98%  '

$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '2.41'
$ws.Range('E47').Value = '  +4.19%  '

$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '26.33'
$ws.Range('E48').Value = '  -1.03%  '

$ws.Range('E49').Value = '  +0.06%  '

$ws.Range('E50').Value = '  -0.46%  '

$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '34.08'
$ws.Range('E51').Value = '  +6.94%  '
